$d = $word.ActiveDocument

# --- Paragraph 1 (title block: date + paper title separated by a line break) ---
$d.Paragraphs.Item(1).Range.Find.Execute("המאמר היומי של מייק - 11.04.25", $false, $false, $false, $false, $false, $true, 1, $false, "המאמר היומי של מייק - 09.04.25", 2) | Out-Null
$d.Paragraphs.Item(1).Range.Find.Execute(" Arithmetic Without Algorithms: Language Models Solve Math with a Bag of Heuristics", $false, $false, $false, $false, $false, $true, 1, $false, "O1-CODER: AN O1 REPLICATION FOR CODING", 2) | Out-Null

# --- Paragraphs 2-8: replace full paragraph body text ---
$d.Paragraphs.Item(2).Range.Find.Execute(" כבר סקרתי בעבר כמה מאמרים על מודלי שפה לחישוב נוסחאות אריתמטיות המכילות פעולות חשבוניות סטנדרטיות כמו פלוס, כפול וכדומה. לדעתי מודלי שפה פחות מיועדים למשימות מהסוג הזה (יש לנו מחשבונים, בפייטון וכאלו) אבל בכל זאת יש מחקרים מעניינים בנושא הזה. ויש סיבה נוספת לבחירת המאמר הזה - הוא נכתב על ידי חוקרים ישראלים ותמיד נהנה לסקור תוצרת מקומית.", $false, $false, $false, $false, $false, $true, 1, $false, "סוף סוף הגעתי לסקור את המאמר הזה שעשה לא מעט רעש בזמנו. המטרה המוצהרת של מחברי המאמר היא לחקות את o1 של OpenAI על משימות קידוד. המאמר השתמש בטכניקות RLHF בשילוב עם שיטת self-play שבה המאמר לומד על הדאטה שהוא עצמו מגנרט. המאמר מתחיל מדאטהסט של של שאלות קוד והשתובות על שאלות אלו (כלומר קוד :)). ", 2) | Out-Null
$d.Paragraphs.Item(3).Range.Find.Execute("אז כאמור המאמר חוקר מה קורה בתוך מודל הטרנספורמר כאשר מודל שפה מקבל משימה אריתמטית. למעשה המחברים מנסים לאתר מה שנקרא נתיב החישובי (circuit) בתוך הטרנספורמר כלומר רכיביו המבצעים בפועל את `"החישובים הנדרשים`" עבור משימה זו. אתם בטח זוכרים שבלוק טרנספורמר מורכב משתי שכבות עיקריות (יש גם שכבות נרמול) שהם מנגנון attention מרובה ראשים או MHA ושכבת MLP המורכבת משתי שכבות לינאריות ואקטיבציה לא לינארית ביניהם. אז הנתיב החישובי מורכב מנוירונים מסוימים בתוך ה-MHA או בתוך ה-MLP.", $false, $false, $false, $false, $false, $true, 1, $false, "הרעיון העיקרי של המאמר מכיל 6 שלבים עיקריים. בשלב הראשון המחברים בונים כלי(המאמר לא מרחיב על זה יותר מדי) לגנרוט טסטים מקיפים עבור שאלת קוד והקוד הנכון עבורה. בהמשך כלי זה(TTG) ישמש לשערוך של ה-reward עבור קוד שנבנה על ידי O1-CODER. ", 2) | Out-Null
$d.Paragraphs.Item(4).Range.Find.Execute("כדי לאתר את הנתיב החישובי, המחברים מבצעים החלפת אקטיבציות (activation patching) של נוירונים בתוך הטרנספורמר המאפשרים לשערך את החשיבות של שכבות MLP וכל ראשי attention בכל מיקום בסדרת קלט (פרומפט אריתמטי). איך עושים זאת? לוקחים פרומפט אריתמטי מסוים (לדוגמא, `"226 − 68 =`"), ופרומפט אקראי שמוביל לתוצאה שונה (למשל, `"21 + 17 =`"). לאחר חישוב של אקטיבציות המודל עבור הפקודה האקראי, מזינים את פרומפט המקורי למודל.", $false, $false, $false, $false, $false, $true, 1, $false, "בשלב השני באמצעות MCTS שזה ראשי תיבות של Markov Chain Tree Search בונים את שרשראות הנמקה (reasoning) עבור הדוגמאות מהדאטהסט. MCTS הוא אלגוריתם לתכנון בקבלת החלטות שמבצע דגימה במרחב המצבים (טוקנים במקרה שלנו) כדי לשערך את הreward הפעולות האפשריות. האלגוריתם בונה עץ חיפוש באופן הדרגתי – בכל צעד הוא בוחר לפתח את הענף(סדרת טוקנים) שנראה הכי מבטיח, תוך איזון בין חקירה של אפשרויות חדשות לבין ניצול של מה שכבר נמצא כמוצלח. כל מסלול בעץ (שרשרת הנמקה הכוללת פתרון) מקבלת תגמול 0 או 1 עם TTG(עובר או לא עובר את כל הטסטים).", 2) | Out-Null
$d.Paragraphs.Item(5).Range.Find.Execute("בשלב זה מתערבים בחישוב (patching) — כלומר, מחליפים את אקטיבציות של שכבת MLP בודדת או ראש attention באקטיבציות שחושבה מראש עבור הפרומפט האקראי. בהמשך בודקים כיצד ההתערבות משפיעה על ההסתברויות של שני הטוקנים של התשובות(עבור הפרופמט המקורי ועבור האקראי) - יש נוסחה שמשערכת השינויים בטוקני התשובות. לאחר מציאת הנתיב החישובי עבור הדוגמאות השונות המאמר משערך את ״נקיונם״ על ידי החלפה של כל האקטיבציות באקטיבציות ממוצעות על פני דאטהסט גדול של פרומפטים אריטמתיים כאשר רק האקטיבציות של הנתיב החישובי נותרו על כנם. המחברים הראו שהחלפה זו כמעט ולא משפיע על הלוגיטים של התשובה הנכונה.", $false, $false, $false, $false, $false, $true, 1, $false, "בשלב השלישי המודל עובר SFT על שרשראות ההנמקה שהובילו לפתרון הנכון (עם ציון 1). בשלב הרביעי מתחילים את אימון ה-self-play בצורה איטרטיבית כאשר דאטהסט האימון מועשר בכל איטרציה עם הדוגמאות הנוצרות על ידי המודל עצמו. בהתחלה מבצעים אימון SFT של המודל על הדאטהטס עם התשובות הנכונות בךבד(פרט לאיטרציה 0) או מבצעים אימון RLHF עם DPO (שזה Direct Preference Optimization) על הזוגות של דוגמאות חיוביות ושליליות. ", 2) | Out-Null
$d.Paragraphs.Item(6).Range.Find.Execute("אחרי מציאת נתיבים חישוביים אלו המחברים ניסו להבין איזה משמעות אריתמטית יש להם. כתוצאה מכך התבררה תמונה די מעניינת. המחברים הראו כי הפעולות של נתיבים אלו הם למעשה יוריסטיקות שונות המאפשרות לפתור את התרגיל. למשל היו נוירונים שמטרתם היא להגיד האם התוצאה נמצאת בתחום [150, 180] או שהתוצאה מתחלקת ב-5. שילוב של שערוכים אלו מאפשר למודל לפתור תרגילים אריתמטיים פשוטים יחסית הלא מערבים מספרים גבוהים מדי. זה די מסביר למה LLMs מתקשים עם פעולות על מספרים גבוהים. ", $false, $false, $false, $false, $false, $true, 1, $false, "לאחר מכן אנו מגנרטים שרשראות הנמקה עם המודל (פרט לתשובה הסופית) ומשתמשים במודל תגמול PRM(שזה Process Reward Model) למתן תגמול לשרשראות הנמקה אלו. אז בונים את התשובה על השאלה משרשרת ההנמקה ויוצרים טסטים לשאלה זו (ידועה לנו התשובה הנכונה לכל שאלה - כנראה השאלות הן חלק מדאטהסט גדול של שאלות פתורות). אחרי זה מחשבים את ה-reward על ידי הרצת טסטים על התשובות שגונרטו על ידי המודל (1 - הטסטים עברו, 0 - לא עברו) ומשלבים אותו עם התגמולים שהתקבלו במהלך הריזונינג (נקרא aggregation function). מאמנים את המודל במטרה למקסם את התגמול הזה (עם שיטת RL כלשהי) - כנראה שיש כאן איזושהי רגולריזציה אבל המאמר לא מרחיב על זה.", 2) | Out-Null
$d.Paragraphs.Item(7).Range.Find.Execute("בנוסף יש כמה מציאות מעניינות. רוב החלקים הבולטים של הנתיבים החישוביים נמצאים בשכבות MLP ולא בראש attention. הדבר המעניין השני הוא העובדה שהמודל ״די מתכנס״ לתשובה הנכונה כבר בשכבות הביניים (ניתן להפיק אותה משם על ידי שכבה לינארית).", $false, $false, $false, $false, $false, $true, 1, $false, "בסוף יוצרים דוגמאות עם המודל אחרי העדכון האחרון ומוסיפים אותם לדאטהסט ומתחילים מחדש את השלב הרביעי (self-play).", 2) | Out-Null
$d.Paragraphs.Item(8).Range.Find.Execute("https://arxiv.org/abs/2410.21272", $false, $false, $false, $false, $false, $true, 1, $false, "מאמר מאוד מעניין…", 2) | Out-Null

# --- Append a new paragraph after the (now-updated) last paragraph, containing the new arXiv link ---
$d.Paragraphs.Item(8).Range.InsertParagraphAfter()
$d.Paragraphs.Item(9).Range.Text = "https://arxiv.org/abs/2412.00154"

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
